$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (E2:T2), keep A2/B2/C2/D2 text as-is (ECs, Trf, Tfr2, FAPs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04202566666666666
$ws.Range("H2").Value = 0.126077
$ws.Range("I2").Value = 0.001003231639737821
$ws.Range("J2").Value = 0.001003231639737821
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.01263189277544444
$ws.Range("R2").Value = 0.113687034979
$ws.Range("S2").Value = 0.001003231639737821
$ws.Range("T2").Value = 0.001003231639737821

# Row 3 becomes the former row4 data (FAPs -> FAPs) but with new TPM-derived specificity values
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 39.70924033333333
$ws.Range("H3").Value = 119.127721
$ws.Range("I3").Value = 0.9479341900351343
$ws.Range("J3").Value = 0.9479341900351345
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3005756666666666
$ws.Range("N3").Value = 0.9017269999999999
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 11.93563138601855
$ws.Range("R3").Value = 107.420682474167
$ws.Range("S3").Value = 0.9479341900351343
$ws.Range("T3").Value = 0.9479341900351345

# Row 4 becomes MuSCs -> FAPs with new TPM-derived values
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.139026333333333
$ws.Range("H4").Value = 6.417078999999999
$ws.Range("I4").Value = 0.05106257832512778
$ws.Range("J4").Value = 0.05106257832512778
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3005756666666666
$ws.Range("N4").Value = 0.9017269999999999
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.6429392661592221
$ws.Range("R4").Value = 5.786453395432999
$ws.Range("S4").Value = 0.05106257832512778
$ws.Range("T4").Value = 0.05106257832512778

# Remove rows 5-7 entirely (they no longer exist in the updated data)
$ws.Range("A5:A7").EntireRow.Delete()
